$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    "61-44=",
    "22-6=",
    "23-18=",
    "6+86=",
    "56+35=",
    "60-58=",
    "90-36=",
    "59+25=",
    "74-67=",
    "3+19=",
    "88+8=",
    "77-29=",
    "8+35=",
    "49+18=",
    "23+28=",
    "90-58=",
    "17+58=",
    "80-35=",
    "44-28=",
    "62-44=",
    "86-48=",
    "46-39=",
    "8+45=",
    "80-77=",
    "49+26=",
    "48+5=",
    "65+18=",
    "65-58=",
    "55-6=",
    "19+46=",
    "55-17=",
    "17+17=",
    "8+5=",
    "10-9=",
    "41-2=",
    "35+59=",
    "36+35=",
    "7+24=",
    "73-8=",
    "39+22=",
    "25+36=",
    "84-75=",
    "76-7=",
    "75+17=",
    "71-47=",
    "76-69=",
    "40-16=",
    "4+18=",
    "47-28=",
    "35+19=",
    "3+39=",
    "64+7=",
    "19+3=",
    "82-56=",
    "35-16=",
    "67-19=",
    "24+19=",
    "19+26=",
    "18+39=",
    "57+19=",
    "69+14=",
    "43+48=",
    "94-38=",
    "59+36=",
    "51-29=",
    "46-7=",
    "54+17=",
    "52-16=",
    "29+38=",
    "54-47=",
    "68+3=",
    "52-44=",
    "61-16=",
    "68+29=",
    "80-6=",
    "74-56=",
    "83-67=",
    "90-41=",
    "39+16=",
    "25+68=",
    "69+9=",
    "12+49=",
    "7+67=",
    "91-83=",
    "16+6=",
    "41-17=",
    "88-9=",
    "61-43=",
    "66-49=",
    "69+3=",
    "18+38=",
    "81-13=",
    "64+7=",
    "15+7=",
    "32-14=",
    "58+8=",
    "50-38=",
    "9+53=",
    "7+57=",
    "13+69="
)

$numCols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
